# Applies the updated crypto price/volume figures for the Mon Apr  8 04:45:03 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds formatted price strings (e.g. "3.417.54") that must stay as text,
# so force a text number format before assigning to avoid Excel auto-converting them to numbers.
$priceCells = @(
    "D2", "D3", "D5", "D6", "D8", "D9", "D10", "D12", "D14", "D15", "D17", "D18", "D21", "D23", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D46", "D47", "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = @{
    'D2' = '69.362.67'
    'E2' = '  +0.27%  '
    'D3' = '3.417.54'
    'E3' = '  +1.08%  '
    'E4' = '  +0.07%  '
    'D5' = '581.62'
    'E5' = '  -0.82%  '
    'D6' = '176.11'
    'E6' = '  -2.33%  '
    'E7' = '  +0.04%  '
    'D8' = '3.413.54'
    'E8' = '  +1.05%  '
    'D9' = '0.590'
    'E9' = '  -0.72%  '
    'D10' = '0.198'
    'E10' = '  +0.30%  '
    'E11' = '  -1.35%  '
    'D12' = '48.50'
    'E12' = '  -0.42%  '
    'E13' = '  -2.34%  '
    'D14' = '691.23'
    'E14' = '  +0.68%  '
    'D15' = '3.964.42'
    'E15' = '  +0.91%  '
    'E16' = '  -0.28%  '
    'D17' = '69.395.02'
    'E17' = '  +0.31%  '
    'D18' = '3.416.57'
    'E18' = '  +1.01%  '
    'E19' = '  +0.80%  '
    'E20' = '  -0.55%  '
    'D21' = '11.31'
    'E21' = '  -0.33%  '
    'E22' = '  -0.99%  '
    'D23' = '5.40'
    'E23' = '  -0.28%  '
    'E24' = '  -1.40%  '
    'D25' = '101.25'
    'E25' = '  -2.59%  '
    'E26' = '  -1.11%  '
    'D27' = '2.65'
    'E27' = '  -2.43%  '
    'E28' = '  -0.97%  '
    'D29' = '33.30'
    'E29' = '  -3.01%  '
    'D30' = '8.71'
    'E30' = '  +0.26%  '
    'D31' = '7.00'
    'E31' = '  +0.38%  '
    'D32' = '575.08'
    'E32' = '  +3.54%  '
    'D33' = '3.68'
    'E33' = '  +0.41%  '
    'D34' = '10.97'
    'E34' = '  -2.01%  '
    'D35' = '58.17'
    'E35' = '  +0.30%  '
    'E36' = '  -3.25%  '
    'D37' = '0.999'
    'E37' = '  -0.07%  '
    'D38' = '3.530.47'
    'E38' = '  -4.64%  '
    'E39' = '  -1.84%  '
    'D40' = '34.70'
    'E40' = '  -0.36%  '
    'D41' = '0.0₃0725'
    'E41' = '  +3.00%  '
    'E42' = '  +0.72%  '
    'E43' = '  -1.02%  '
    'E44' = '  -2.36%  '
    'E45' = '  -0.67%  '
    'D46' = '1.44'
    'E46' = '  +3.97%  '
    'D47' = '2.63'
    'E47' = '  -0.32%  '
    'E48' = '  -1.10%  '
    'E49' = '  -0.13%  '
    'D50' = '132.55'
    'E50' = '  +0.04%  '
    'E51' = '  +1.79%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
